$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (column G) values per row, per diff
$kValues = @{
    2 = 0
    3 = 2
    4 = 0
    5 = 1
    6 = 0
    7 = 2
    8 = 1
    9 = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 1
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 2
    29 = 1
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 0
    35 = 3
    36 = 0
    37 = 2
    38 = 2
    39 = 0
    40 = 1
    41 = 2
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 2
    48 = 1
    49 = 1
    50 = 1
    51 = 3
    52 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
